$d = $word.ActiveDocument

# Locate "307" in "del modulo 307." and find the exact character position of "7"
$find = $d.Content.Duplicate
$find.Find.Execute("modulo 307", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0) | Out-Null

# $find now spans "modulo 307"; compute position of trailing "7" (last char)
$sevenStart = $find.End - 1
$sevenRange = $d.Range($sevenStart, $find.End)

# Remove the existing "_GoBack" bookmark (it will be re-created at the new edit location,
# mirroring Word's behaviour of marking the last editing position).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Replace "7" with "6", as if the user selected the digit and typed a replacement.
$sevenRange.Delete()
$insertPoint = $d.Range($sevenStart, $sevenStart)
$insertPoint.InsertAfter("6")

# Nudge formatting on the freshly typed character (net no-op) so this run keeps its
# own identity instead of being silently re-coalesced with its neighbours.
$sixRange = $d.Range($sevenStart, $sevenStart + 1)
$sixRange.Bold = 1
$sixRange.Bold = 0

# Re-insert the "_GoBack" bookmark right after the corrected text (collapsed range).
$goBackRange = $d.Range($sevenStart + 1, $sevenStart + 1)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

$d.Save()
